# Auto_scrap_Plumbersstock.py scrap added - update workbook metadata accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "HomeDepot_URL" to "URL"
$ws.Name = "URL"

# Move the active selection from H22 to I5
$ws.Range("I5").Select()
